$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.032.52'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.23%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.518.11'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '609.68'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.93'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.71%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.518.33'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.53%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.478'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.65%  '
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '8.09'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +7.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.423'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.70%  '
$ws.Range('E13').Value = '  +1.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.113.66'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '31.89'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.522.53'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.068.44'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.23%  '
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.81'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +8.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.46'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.42'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '437.85'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.76%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.610'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.70'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.656.38'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.38%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  -4.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.78'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.29'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -4.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.53'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('E31').Value = '  -2.51%  '
$ws.Range('E32').Value = '  -1.69%  '
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.61'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.98'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.06%  '
$ws.Range('E36').Value = '  -2.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.08'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '175.93'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.05'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -12.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.895'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.16'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -7.16%  '
$ws.Range('E47').Value = '  -2.70%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.48'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.77%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.46'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.997'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.53%  '
